$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the WhatsApp Number value in E2
$ws.Range("E2").Value = 8145312848

# Move the active selection to E2 (matches the saved view state in the diff)
$ws.Range("E2").Select() | Out-Null
